# Auto-generated edit script: "Generate Report for Handoff"
# Re-sorts the three status rows alphabetically by file name and updates the
# status/date fields for the 7f9414b6 file to reflect a fresh handoff, then
# rebuilds the hyperlinks to match the new row layout.

$wb = $excel.ActiveWorkbook

# ---- Overview ----
$ws = $wb.Worksheets.Item("Overview")

# Update cell values for the re-sorted rows
$ws.Range("A2").Value2 = "ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "2016-03-24 10:41:38"
$ws.Range("A3").Value2 = "ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("D3").Value2 = "2016-03-24 10:41:38"
$ws.Range("A4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "Ready for handoff"
$ws.Range("D4").Value2 = "2016-03-24 10:45:33"

# Rebuild hyperlinks to match the new row order (this also clears the old ones)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md", "", "", "ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md", "", "", "ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md")

# ---- zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values for the re-sorted rows
$ws.Range("A2").Value2 = "ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md"
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf"
$ws.Range("E2").Value2 = "2016-03-24 10:41:29"
$ws.Range("F2").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.md"
$ws.Range("G2").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf"
$ws.Range("H2").Value2 = "2016-03-24 10:42:10"
$ws.Range("J2").Value2 = "Include"
$ws.Range("A3").Value2 = "ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md"
$ws.Range("B3").Value2 = ".md"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("D3").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf"
$ws.Range("E3").Value2 = "2016-03-24 10:41:29"
$ws.Range("F3").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.md"
$ws.Range("G3").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf"
$ws.Range("H3").Value2 = "2016-03-24 10:42:10"
$ws.Range("J3").Value2 = "Include"
$ws.Range("A4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md"
$ws.Range("B4").Value2 = ".md"
$ws.Range("C4").Value2 = "Ready for handoff"
$ws.Range("D4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.zh-cn.xlf"
$ws.Range("E4").Value2 = "2016-03-24 10:45:23"
$ws.Range("F4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md"
$ws.Range("G4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.zh-cn.xlf"
$ws.Range("H4").Value2 = "2016-03-24 10:44:25"
$ws.Range("J4").Value2 = "Include"

# Rebuild hyperlinks to match the new row order (this also clears the old ones)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md", "", "", "ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b0dadece75fc4157b6cc34eb2f3cd5e8c624a31/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.zh-cn.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d793c16a862437f073c896b53ac6b096e0e94106/e2e/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dbdbcdb77d65450e55ce126221df9da51328a967/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.zh-cn.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md", "", "", "ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3855c023bd6016ce0248ee16958e52ab4be1bd79/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/b60407fda23992a3cf75b09ca88f2a6dce583a24/e2e/6bda3337-91f3-4f06-bf98-fa00d01369a2.md", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/05cfa108d1585967a9d35d84c462033b827639f2/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3855c023bd6016ce0248ee16958e52ab4be1bd79/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/b60407fda23992a3cf75b09ca88f2a6dce583a24/e2e/6bda3337-91f3-4f06-bf98-fa00d01369a2.md", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/05cfa108d1585967a9d35d84c462033b827639f2/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.zh-cn.xlf", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.zh-cn.xlf")

# ---- de-de ----
$ws = $wb.Worksheets.Item("de-de")

# Update cell values for the re-sorted rows
$ws.Range("A2").Value2 = "ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md"
$ws.Range("B2").Value2 = ".md"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf"
$ws.Range("E2").Value2 = "2016-03-24 10:41:38"
$ws.Range("F2").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.md"
$ws.Range("G2").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf"
$ws.Range("H2").Value2 = "2016-03-24 10:42:25"
$ws.Range("J2").Value2 = "Include"
$ws.Range("A3").Value2 = "ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md"
$ws.Range("B3").Value2 = ".md"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"
$ws.Range("D3").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf"
$ws.Range("E3").Value2 = "2016-03-24 10:41:38"
$ws.Range("F3").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.md"
$ws.Range("G3").Value2 = "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf"
$ws.Range("H3").Value2 = "2016-03-24 10:42:25"
$ws.Range("J3").Value2 = "Include"
$ws.Range("A4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md"
$ws.Range("B4").Value2 = ".md"
$ws.Range("C4").Value2 = "Ready for handoff"
$ws.Range("D4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.de-de.xlf"
$ws.Range("E4").Value2 = "2016-03-24 10:45:33"
$ws.Range("F4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md"
$ws.Range("G4").Value2 = "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.de-de.xlf"
$ws.Range("H4").Value2 = "2016-03-24 10:44:39"
$ws.Range("J4").Value2 = "Include"

# Rebuild hyperlinks to match the new row order (this also clears the old ones)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md", "", "", "ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2fb5f702a648ce0ebce145b999a32ac00ffb1474/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.de-de.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/957fd60aa76a2c3089fae931cd96fc8323f1525f/e2e/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b1c320d94165f93a97f8e8a0dd29be4c2232e19c/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.de-de.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/ffffa84066c0-4ad8-4e0e-a2a2-6dbf8b12e950.md", "", "", "ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3d79365a3c67b449755f407a72a725e6c29f7d7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/098074e0217346ce75db31a83e76b291db4e6078/e2e/6bda3337-91f3-4f06-bf98-fa00d01369a2.md", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.md")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/588819832a0e8a8cc1c3163ff9ddd7a79a6a6482/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf", "", "", "6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e66c7bc00d375c333839c52ef0d4d02ef7602f5e/e2e/ffffff6b5652de-99b3-409d-8f7d-209d2a47b918.md", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3d79365a3c67b449755f407a72a725e6c29f7d7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/098074e0217346ce75db31a83e76b291db4e6078/e2e/6bda3337-91f3-4f06-bf98-fa00d01369a2.md", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.md")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/588819832a0e8a8cc1c3163ff9ddd7a79a6a6482/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/6bda3337-91f3-4f06-bf98-fa00d01369a2.80fafdefc26409659540283adafd43685004e993.de-de.xlf", "", "", "7f9414b6-db4a-46bf-9691-d99c0d3b66f9.62217c803f10f18cf603cde202caa09c7a1f5f24.de-de.xlf")
